$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 158.63637
$ws.Range("I11").Value = 158.63637
$ws.Range("K11").Value = 158.63637
$ws.Range("M11").Value = -18.63637
$ws.Range("H17").Value = 1628.1111
$ws.Range("J17").Value = 1706.625
$ws.Range("L17").Value = 5119.875
$ws.Range("N17").Value = -5455.875
$ws.Range("H32").Value = 4464.95
$ws.Range("I32").Value = 6576.909
$ws.Range("J32").Value = 1883.6666
$ws.Range("K32").Value = 6576.909
$ws.Range("L32").Value = 1883.6666
$ws.Range("M32").Value = -6250.909
$ws.Range("N32").Value = -2535.6666
$ws.Range("H51").Value = 9561.214
$ws.Range("I51").Value = 4744.25
$ws.Range("J51").Value = 11488
$ws.Range("K51").Value = 4744.25
$ws.Range("L51").Value = 11488
$ws.Range("M51").Value = -4260.25
$ws.Range("N51").Value = -12456
$ws.Range("H62").Value = 4888.892
$ws.Range("I62").Value = 4764.4443
$ws.Range("J62").Value = 5224.9
$ws.Range("K62").Value = 4764.4443
$ws.Range("L62").Value = 5224.9
$ws.Range("M62").Value = -4140.4443
$ws.Range("N62").Value = -6472.9
$ws.Range("H65").Value = 4888.892
$ws.Range("I65").Value = 4764.4443
$ws.Range("J65").Value = 5224.9
$ws.Range("K65").Value = 23822.2215
$ws.Range("L65").Value = 26124.5
$ws.Range("M65").Value = -20702.2215
$ws.Range("N65").Value = -32364.5
$ws.Range("H88").Value = 5306.846
$ws.Range("I88").Value = 1923.75
$ws.Range("J88").Value = 10719.8
$ws.Range("K88").Value = 1923.75
$ws.Range("L88").Value = 10719.8
$ws.Range("M88").Value = -1517.75
$ws.Range("N88").Value = -11531.8
$ws.Range("H91").Value = 5306.846
$ws.Range("I91").Value = 1923.75
$ws.Range("J91").Value = 10719.8
$ws.Range("K91").Value = 1923.75
$ws.Range("L91").Value = 10719.8
$ws.Range("M91").Value = -519.75
$ws.Range("N91").Value = -13527.8
$ws.Range("H116").Value = 14574.226
$ws.Range("I116").Value = 16996.13
$ws.Range("K116").Value = 16996.13
$ws.Range("M116").Value = -13554.13
$ws.Range("H132").Value = 21015.889
$ws.Range("I132").Value = 23976.516
$ws.Range("K132").Value = 71929.548
$ws.Range("M132").Value = -69399.548
$ws.Range("H138").Value = 32858.426
$ws.Range("I138").Value = 2173.45
$ws.Range("J138").Value = 80066.08
$ws.Range("K138").Value = 6520.349999999999
$ws.Range("L138").Value = 240198.24
$ws.Range("M138").Value = -1380.349999999999
$ws.Range("N138").Value = -250478.24

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20663.46
$ws.Range("I32").Value = 22268.854
$ws.Range("K32").Value = 22268.854
$ws.Range("M32").Value = -21981.854
$ws.Range("H39").Value = 9229
$ws.Range("I39").Value = 9229
$ws.Range("K39").Value = 9229
$ws.Range("M39").Value = -8709
$ws.Range("H45").Value = 4380.3076
$ws.Range("J45").Value = 6205.7144
$ws.Range("L45").Value = 6205.7144
$ws.Range("N45").Value = -6959.7144
$ws.Range("H74").Value = 382985.44
$ws.Range("I74").Value = 600830.6
$ws.Range("K74").Value = 600830.6
$ws.Range("M74").Value = -599956.6
$ws.Range("H77").Value = 382985.44
$ws.Range("I77").Value = 600830.6
$ws.Range("K77").Value = 3004153
$ws.Range("M77").Value = -2999785
$ws.Range("H97").Value = 1696.6842
$ws.Range("I97").Value = 1295.75
$ws.Range("K97").Value = 1295.75
$ws.Range("M97").Value = -799.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 2195.5
$ws.Range("J20").Value = 1892.1111
$ws.Range("K20").Value = 2195.5
$ws.Range("L20").Value = 1892.1111
$ws.Range("M20").Value = -1948.5
$ws.Range("N20").Value = -2386.1111
$ws.Range("H64").Value = 1323.7273
$ws.Range("J64").Value = 1114.2222
$ws.Range("L64").Value = 1114.2222
$ws.Range("N64").Value = -1564.2222
$ws.Range("H67").Value = 1323.7273
$ws.Range("J67").Value = 1114.2222
$ws.Range("L67").Value = 1114.2222
$ws.Range("N67").Value = -2674.2222
$ws.Range("H107").Value = 2187.1765
$ws.Range("I107").Value = 2019.6522
$ws.Range("J107").Value = 2537.4546
$ws.Range("K107").Value = 2019.6522
$ws.Range("L107").Value = 2537.4546
$ws.Range("M107").Value = -99.65219999999999
$ws.Range("N107").Value = -6377.4546

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1700
$ws.Range("I16").Value = 1700
$ws.Range("K16").Value = 1700
$ws.Range("M16").Value = -1413
$ws.Range("H58").Value = 1168.96
$ws.Range("I58").Value = 932.75
$ws.Range("K58").Value = 932.75
$ws.Range("M58").Value = -729.75
$ws.Range("H99").Value = 6985
$ws.Range("I99").Value = 5242.6665
$ws.Range("K99").Value = 5242.6665
$ws.Range("M99").Value = -3744.6665
$ws.Range("H113").Value = 1700
$ws.Range("I113").Value = 1700
$ws.Range("K113").Value = 1700
$ws.Range("M113").Value = 470
$ws.Range("H126").Value = 6985
$ws.Range("I126").Value = 5242.6665
$ws.Range("K126").Value = 15727.9995
$ws.Range("M126").Value = -13257.9995
$ws.Range("H132").Value = 113871.336
$ws.Range("I132").Value = 250752.25
$ws.Range("K132").Value = 752256.75
$ws.Range("M132").Value = -749726.75
$ws.Range("H136").Value = 1168.96
$ws.Range("I136").Value = 932.75
$ws.Range("K136").Value = 2798.25
$ws.Range("M136").Value = -248.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2629.9
$ws.Range("J34").Value = 2888.7778
$ws.Range("L34").Value = 8666.3334
$ws.Range("N34").Value = -8834.3334
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 12000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -12354

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4793.4375
$ws.Range("I70").Value = 4615.25
$ws.Range("J70").Value = 4971.625
$ws.Range("K70").Value = 4615.25
$ws.Range("L70").Value = 4971.625
$ws.Range("M70").Value = -4345.25
$ws.Range("N70").Value = -5511.625
$ws.Range("H73").Value = 4793.4375
$ws.Range("I73").Value = 4615.25
$ws.Range("J73").Value = 4971.625
$ws.Range("K73").Value = 4615.25
$ws.Range("L73").Value = 4971.625
$ws.Range("M73").Value = -3679.25
$ws.Range("N73").Value = -6843.625
$ws.Range("H97").Value = 1223.591
$ws.Range("I97").Value = 1173.8889
$ws.Range("K97").Value = 1173.8889
$ws.Range("M97").Value = -677.8888999999999
$ws.Range("H122").Value = 3747.2258
$ws.Range("I122").Value = 3350.6956
$ws.Range("J122").Value = 4887.25
$ws.Range("K122").Value = 10052.0868
$ws.Range("L122").Value = 14661.75
$ws.Range("M122").Value = -7602.086800000001
$ws.Range("N122").Value = -19561.75
$ws.Range("H126").Value = 2444.3333
$ws.Range("I126").Value = 1999.875
$ws.Range("K126").Value = 5999.625
$ws.Range("M126").Value = -3529.625

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5030000
$ws.Range("J2").Value = 5030000
$ws.Range("L2").Value = 5030000
$ws.Range("N2").Value = -5030224
$ws.Range("H7").Value = 4542.375
$ws.Range("J7").Value = 3600
$ws.Range("L7").Value = 3600
$ws.Range("N7").Value = -3824
$ws.Range("H40").Value = 3766.1875
$ws.Range("I40").Value = 3447.0715
$ws.Range("K40").Value = 3447.0715
$ws.Range("M40").Value = -3311.0715
$ws.Range("H46").Value = 5413.4443
$ws.Range("I46").Value = 1750.4
$ws.Range("K46").Value = 1750.4
$ws.Range("M46").Value = -1562.4
$ws.Range("H54").Value = 34992.5
$ws.Range("I54").Value = 34991
$ws.Range("J54").Value = 34994
$ws.Range("K54").Value = 34991
$ws.Range("L54").Value = 34994
$ws.Range("M54").Value = -34347
$ws.Range("N54").Value = -36282
$ws.Range("H68").Value = 3833.1177
$ws.Range("I68").Value = 3030.111
$ws.Range("K68").Value = 3030.111
$ws.Range("M68").Value = -2281.111
$ws.Range("H71").Value = 3833.1177
$ws.Range("I71").Value = 3030.111
$ws.Range("K71").Value = 15150.555
$ws.Range("M71").Value = -11406.555
$ws.Range("H93").Value = 1258.762
$ws.Range("I93").Value = 1145.0588
$ws.Range("K93").Value = 1145.0588
$ws.Range("M93").Value = 102.9412
$ws.Range("H126").Value = 4542.375
$ws.Range("J126").Value = 3600
$ws.Range("L126").Value = 10800
$ws.Range("N126").Value = -15740
$ws.Range("H132").Value = 4995.2856
$ws.Range("I132").Value = 4995.2856
$ws.Range("K132").Value = 14985.8568
$ws.Range("M132").Value = -12455.8568
$ws.Range("H136").Value = 3756
$ws.Range("I136").Value = 3155
$ws.Range("K136").Value = 9465
$ws.Range("M136").Value = -6915

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1240.3846
$ws.Range("J107").Value = 929.3333
$ws.Range("L107").Value = 2787.9999
$ws.Range("N107").Value = -6627.9999
$ws.Range("H132").Value = 20582.297
$ws.Range("I132").Value = 23611.129
$ws.Range("K132").Value = 70833.387
$ws.Range("M132").Value = -68303.387
